# Replace the trailing "nan" label (produced from a pandas melt/crosstab
# "Total" row whose group label was missing) with the literal word "Total"
# for the five distinct prefixes that appear in this workbook.
$wb = $excel.ActiveWorkbook

$replacements = @{
    "None nan"                        = "None Total"
    "Current nan"                     = "Current Total"
    "Cash flow hedge reserve nan"     = "Cash flow hedge reserve Total"
    "Reserve in joint venture nan"    = "Reserve in joint venture Total"
    "Other reserve nan"               = "Other reserve Total"
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    if ($used -eq $null) { continue }

    $rows = $used.Rows.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        $cell = $ws.Cells.Item($startRow + $r, $startCol)
        $val = $cell.Value2
        if ($val -ne $null -and $replacements.ContainsKey([string]$val)) {
            $cell.Value = $replacements[[string]$val]
        }
    }
}
